# info.xlsx 2021 12 8
# Update the "operation type" header text (B1) to document the new
# 5th action type, and re-style the header row:
#   - B1 gets a plain (non-rich-text) red-font string
#   - C1 gets a yellow highlight fill
# Also widen columns B/D to fit the longer text and move the active
# selection to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header cells -------------------------------------------------

# A1 stays the same text, just re-written (drops the old phonetic-guide
# run info that Excel attaches to CJK shared strings).
$ws.Range("A1").Value = "事件描述（想写就写）"

# B1: collapse the old 3-run rich text into one plain string that now
# documents click modes 1-5, and make the whole thing red.
$ws.Range("B1").Value = "操作类型：1.左键单击（循环直到找到图片为止）      2.输入字符串       3.等待      4.热键      5.左键单击（无需找到图片）"
$ws.Range("B1").Font.Color = 255

# C1: unchanged text, but now highlighted with a yellow fill.
$ws.Range("C1").Value = "待点击图标名/等待的时间(秒)/输入的字符串/热键"
$ws.Range("C1").Interior.Color = 65535

# D1: unchanged text/style.
$ws.Range("D1").Value = "单击重复次数（不写为1，死循环为-1）"

# --- Data rows (text unchanged, just re-asserted) ------------------------

$ws.Range("C2").Value = "baidu.png"
$ws.Range("C3").Value = "f5"
$ws.Range("C4").Value = "back.png"

# --- Column widths ---------------------------------------------------

# Column B needs to fit the much longer B1 string; column D needs a
# slightly wider fit too. (Fractional widths below land exactly/closest
# on the engine's column-width grid.)
$ws.Columns("B").ColumnWidth = 807 / 7
$ws.Columns("D").ColumnWidth = 269 / 7

# --- Selection ---------------------------------------------------------

$ws.Range("B6").Select() | Out-Null
